$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2318
$ws.Range("F3").Value = 390
$ws.Range("F5").Value = 311
$ws.Range("F6").Value = 311
$ws.Range("F7").Value = 506
$ws.Range("F9").Value = 739
$ws.Range("F11").Value = 735
$ws.Range("F12").Value = 380
$ws.Range("F13").Value = 78
$ws.Range("F16").Value = 999
$ws.Range("F17").Value = 18384
$ws.Range("G17").Value = "暂时售罄"
$ws.Range("F18").Value = 440
$ws.Range("F20").Value = 191
$ws.Range("F21").Value = 272
$ws.Range("F22").Value = 163
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 167
$ws.Range("F28").Value = 309
$ws.Range("F29").Value = 126

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 69
$ws.Range("F7").Value = 196
$ws.Range("F8").Value = 215
$ws.Range("F9").Value = 3338
$ws.Range("F11").Value = 64
$ws.Range("F15").Value = 115
$ws.Range("F17").Value = 2786

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 261
$ws.Range("F3").Value = 75
$ws.Range("F4").Value = 526
$ws.Range("F5").Value = 195

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 261
$ws.Range("F3").Value = 75
$ws.Range("F4").Value = 69
$ws.Range("F6").Value = 2318
$ws.Range("F7").Value = 526
$ws.Range("F8").Value = 390
$ws.Range("F10").Value = 311
$ws.Range("F11").Value = 311
$ws.Range("F12").Value = 506
$ws.Range("F17").Value = 196
$ws.Range("F18").Value = 195
$ws.Range("F19").Value = 739
$ws.Range("F21").Value = 735
$ws.Range("F22").Value = 380
$ws.Range("F23").Value = 78
$ws.Range("F26").Value = 999
$ws.Range("F27").Value = 18384
$ws.Range("F28").Value = 215
$ws.Range("F29").Value = 3338
$ws.Range("F31").Value = 64
$ws.Range("F33").Value = 440
$ws.Range("F35").Value = 191
$ws.Range("F38").Value = 272
$ws.Range("F39").Value = 163
$ws.Range("F42").Value = 1
$ws.Range("F43").Value = 115
$ws.Range("F45").Value = 167
$ws.Range("F47").Value = 309
$ws.Range("F48").Value = 126
$ws.Range("F49").Value = 2786
